$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 392
$ws1.Range("F4").Value = 1632
$ws1.Range("F7").Value = 417
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 0
$ws1.Range("F10").Value = 496

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 392
$ws4.Range("F5").Value = 0
$ws4.Range("F6").Value = 23
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 0
